$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension -> measure for continente/area; aragon dimension -> sdmx refArea
$ws.Range("A2").Value = "iaest-measure:continente"
$ws.Range("B2").Value = "iaest-measure:area"
$ws.Range("H2").Value = "sdmx-dimension:refArea"

# Row 3: dim -> medida for continente/area columns
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "medida"

# Row 4: skos:Concept -> xsd:int for continente/area columns; aragon URI type
$ws.Range("A4").Value = "xsd:int"
$ws.Range("B4").Value = "xsd:int"
$ws.Range("H4").Value = "URI-Comunidad"

# Row 5 (mapping file references) is no longer needed - remove it entirely
$ws.Rows.Item(5).Delete()
